{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Two textual changes (per the target diff):\n//   1. \"App\" -> \"AppConfig\"   (append \"Config\" right after the word \"App\")\n//   2. \"SMSs\" -> \"SMS\\u2019s\" (turn the mis-typed plural \"SMSs\" into the\n//      possessive/contraction \"SMS's\", using a right single quote)\n\nconst body = context.document.body;\n\n// --- Change 1: \"App\" -> \"AppConfig\" -------------------------------------\nconst appResults = body.search(\"App\", { matchCase: true, matchWholeWord: true });\nappResults.load(\"items\");\nawait context.sync();\n\nif (appResults.items.length > 0) {\n  // Insert \"Config\" immediately after the matched \"App\" range so the\n  // visible word becomes \"AppConfig\".\n  appResults.items[0].insertText(\"Config\", Word.InsertLocation.end);\n  await context.sync();\n}\n\n// --- Change 2: \"SMSs\" -> \"SMS\u2019s\" -----------------------------------------\nconst smsResults = body.search(\"SMSs\", { matchCase: true, matchWholeWord: true });\nsmsResults.load(\"items\");\nawait context.sync();\n\nif (smsResults.items.length > 0) {\n  // Replace the whole matched word with the corrected text.\n  smsResults.items[0].insertText(\"SMS\\u2019s\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d below.\n#\n# Two textual changes (per the target diff):\n#   1. \"App\" -> \"AppConfig\"   (append \"Config\" right after the word \"App\")\n#   2. \"SMSs\" -> \"SMS's\"      (turn the mis-typed plural \"SMSs\" into the\n#      possessive/contraction \"SMS's\", using a right single quote U+2019)\n\n$d = $word.ActiveDocument\n\n# --- Change 1: \"App\" -> \"AppConfig\" --------------------------------------\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"App\"\n$rng.Find.MatchWholeWord = $true\n$rng.Find.MatchCase = $true\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0\n$found = $rng.Find.Execute()\nif ($found) {\n    # Collapse to the end of the found \"App\" and append \"Config\" right\n    # after it, so the visible word becomes \"AppConfig\".\n    $rng.Collapse(0)\n    $rng.InsertAfter(\"Config\")\n}\n\n# --- Change 2: \"SMSs\" -> \"SMS's\" ------------------------------------------\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Text = \"SMSs\"\n$rng2.Find.MatchWholeWord = $true\n$rng2.Find.MatchCase = $true\n$rng2.Find.Forward = $true\n$rng2.Find.Wrap = 0\n$found2 = $rng2.Find.Execute()\nif ($found2) {\n    $rng2.Text = \"SMS\" + [char]0x2019 + \"s\"\n}\n"}
